$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at 123. This shifts the old "last row" (123, with the
#    thicker bottom border) down to 124, and shifts the signature block
#    (old rows 128-129) down to rows 129-130.
$ws.Rows("123:123").Insert()

# 2. The newly inserted row 123 is blank; give it the same formatting as the
#    regular data rows above it (row 122) via a formats-only paste.
$ws.Range("B122:J122").Copy()
$ws.Range("B123:J123").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3. Fill in the repeating "Tipo Doc Trabajador / N Doc Trabajador / Nombre
#    Trabajador" / "Salario Basico" columns for the new row 123 (same worker
#    and salary as every other row).
$ws.Range("B123").Value = "CC"
$ws.Range("C123").Value = "1047471451"
$ws.Range("D123").Value = "JEAN CARLOS MORALES TORRES"
$ws.Range("G123").Value = 781242

# 4. Re-populate "Periodo Mora" (E) and "Valor Mora" (F) for every data row
#    (16-124): a new period (2507) was added and the list is now sorted from
#    the newest period down to the oldest.
$arr = New-Object 'object[,]' 109,2
$arr[0,0] = "2507"; $arr[0,1] = 31249
$arr[1,0] = "2506"; $arr[1,1] = 31249
$arr[2,0] = "2505"; $arr[2,1] = 31249
$arr[3,0] = "2504"; $arr[3,1] = 31249
$arr[4,0] = "2503"; $arr[4,1] = 31249
$arr[5,0] = "2502"; $arr[5,1] = 31249
$arr[6,0] = "2501"; $arr[6,1] = 31249
$arr[7,0] = "2412"; $arr[7,1] = 31249
$arr[8,0] = "2411"; $arr[8,1] = 31249
$arr[9,0] = "2410"; $arr[9,1] = 31249
$arr[10,0] = "2409"; $arr[10,1] = 31249
$arr[11,0] = "2408"; $arr[11,1] = 31249
$arr[12,0] = "2407"; $arr[12,1] = 31249
$arr[13,0] = "2406"; $arr[13,1] = 31249
$arr[14,0] = "2405"; $arr[14,1] = 31249
$arr[15,0] = "2404"; $arr[15,1] = 31249
$arr[16,0] = "2403"; $arr[16,1] = 31249
$arr[17,0] = "2402"; $arr[17,1] = 31249
$arr[18,0] = "2401"; $arr[18,1] = 31249
$arr[19,0] = "2312"; $arr[19,1] = 31249
$arr[20,0] = "2311"; $arr[20,1] = 31249
$arr[21,0] = "2310"; $arr[21,1] = 31249
$arr[22,0] = "2309"; $arr[22,1] = 31249
$arr[23,0] = "2308"; $arr[23,1] = 31249
$arr[24,0] = "2307"; $arr[24,1] = 31249
$arr[25,0] = "2306"; $arr[25,1] = 31249
$arr[26,0] = "2305"; $arr[26,1] = 31249
$arr[27,0] = "2304"; $arr[27,1] = 31249
$arr[28,0] = "2303"; $arr[28,1] = 31249
$arr[29,0] = "2302"; $arr[29,1] = 31249
$arr[30,0] = "2301"; $arr[30,1] = 31249
$arr[31,0] = "2212"; $arr[31,1] = 31249
$arr[32,0] = "2211"; $arr[32,1] = 31249
$arr[33,0] = "2210"; $arr[33,1] = 31249
$arr[34,0] = "2209"; $arr[34,1] = 31249
$arr[35,0] = "2208"; $arr[35,1] = 31249
$arr[36,0] = "2207"; $arr[36,1] = 31249
$arr[37,0] = "2206"; $arr[37,1] = 31249
$arr[38,0] = "2205"; $arr[38,1] = 31249
$arr[39,0] = "2204"; $arr[39,1] = 31249
$arr[40,0] = "2203"; $arr[40,1] = 31249
$arr[41,0] = "2202"; $arr[41,1] = 31249
$arr[42,0] = "2201"; $arr[42,1] = 31249
$arr[43,0] = "2112"; $arr[43,1] = 31249
$arr[44,0] = "2111"; $arr[44,1] = 31249
$arr[45,0] = "2110"; $arr[45,1] = 31249
$arr[46,0] = "2109"; $arr[46,1] = 31249
$arr[47,0] = "2108"; $arr[47,1] = 31249
$arr[48,0] = "2107"; $arr[48,1] = 31249
$arr[49,0] = "2106"; $arr[49,1] = 31249
$arr[50,0] = "2105"; $arr[50,1] = 31249
$arr[51,0] = "2104"; $arr[51,1] = 31249
$arr[52,0] = "2103"; $arr[52,1] = 31249
$arr[53,0] = "2102"; $arr[53,1] = 31249
$arr[54,0] = "2101"; $arr[54,1] = 31249
$arr[55,0] = "2012"; $arr[55,1] = 31249
$arr[56,0] = "2011"; $arr[56,1] = 31249
$arr[57,0] = "2010"; $arr[57,1] = 31249
$arr[58,0] = "2009"; $arr[58,1] = 31249
$arr[59,0] = "2008"; $arr[59,1] = 31249
$arr[60,0] = "2007"; $arr[60,1] = 31249
$arr[61,0] = "2006"; $arr[61,1] = 31249
$arr[62,0] = "2005"; $arr[62,1] = 31249
$arr[63,0] = "2004"; $arr[63,1] = 31249
$arr[64,0] = "2003"; $arr[64,1] = 31249
$arr[65,0] = "2002"; $arr[65,1] = 31249
$arr[66,0] = "2001"; $arr[66,1] = 31249
$arr[67,0] = "1912"; $arr[67,1] = 31249
$arr[68,0] = "1911"; $arr[68,1] = 31249
$arr[69,0] = "1910"; $arr[69,1] = 31249
$arr[70,0] = "1909"; $arr[70,1] = 31249
$arr[71,0] = "1908"; $arr[71,1] = 31249
$arr[72,0] = "1907"; $arr[72,1] = 31249
$arr[73,0] = "1906"; $arr[73,1] = 31249
$arr[74,0] = "1905"; $arr[74,1] = 31249
$arr[75,0] = "1904"; $arr[75,1] = 31249
$arr[76,0] = "1903"; $arr[76,1] = 31249
$arr[77,0] = "1902"; $arr[77,1] = 31249
$arr[78,0] = "1901"; $arr[78,1] = 31249
$arr[79,0] = "1812"; $arr[79,1] = 31249
$arr[80,0] = "1811"; $arr[80,1] = 31249
$arr[81,0] = "1810"; $arr[81,1] = 31249
$arr[82,0] = "1809"; $arr[82,1] = 31249
$arr[83,0] = "1808"; $arr[83,1] = 29896
$arr[84,0] = "1807"; $arr[84,1] = 29896
$arr[85,0] = "1806"; $arr[85,1] = 29896
$arr[86,0] = "1805"; $arr[86,1] = 29896
$arr[87,0] = "1804"; $arr[87,1] = 29896
$arr[88,0] = "1803"; $arr[88,1] = 29896
$arr[89,0] = "1802"; $arr[89,1] = 29896
$arr[90,0] = "1801"; $arr[90,1] = 29896
$arr[91,0] = "1712"; $arr[91,1] = 29896
$arr[92,0] = "1711"; $arr[92,1] = 29896
$arr[93,0] = "1710"; $arr[93,1] = 29896
$arr[94,0] = "1709"; $arr[94,1] = 29896
$arr[95,0] = "1708"; $arr[95,1] = 29896
$arr[96,0] = "1707"; $arr[96,1] = 29896
$arr[97,0] = "1706"; $arr[97,1] = 29896
$arr[98,0] = "1705"; $arr[98,1] = 29896
$arr[99,0] = "1704"; $arr[99,1] = 29896
$arr[100,0] = "1703"; $arr[100,1] = 29896
$arr[101,0] = "1702"; $arr[101,1] = 29896
$arr[102,0] = "1701"; $arr[102,1] = 29896
$arr[103,0] = "1612"; $arr[103,1] = 29896
$arr[104,0] = "1611"; $arr[104,1] = 29896
$arr[105,0] = "1610"; $arr[105,1] = 29896
$arr[106,0] = "1609"; $arr[106,1] = 29896
$arr[107,0] = "1608"; $arr[107,1] = 29896
$arr[108,0] = "1607"; $arr[108,1] = 29896

$ws.Range("E16:F124").Value = $arr

# 5. Update the summary figures at the top of the sheet.
$ws.Range("E11").Value = 3370963
$ws.Range("F13").Value = 109
